$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

$ws2.Range("A11:D11").Copy()
$ws2.Range("A12:D12").PasteSpecial(-4122)

$ws2.Range("A12").Value = "v2.0"
$ws2.Range("B12").Value = "Mahmoud Abdelmageed"
$ws2.Range("C12").Value = "Reviewed publish video"
$ws2.Range("D12").Value = 45785
